# Weekly update: a new price-report week is inserted for this product.
# Two new rows (quality grades "Primera" and "Especial") are inserted right
# after the current most-recent week (row 43), pushing all the older weeks
# down by two rows. The two new rows are then populated with this week's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 44 (old rows 44:51 shift down to 46:53).
$ws.Rows("44:45").Insert()

# New row 44: "Primera" quality for 2023-06-29 (serial 45106)
$ws.Cells.Item(44, 1).Value  = 7
$ws.Cells.Item(44, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(44, 3).Value  = "Ñuble"
$ws.Cells.Item(44, 4).Value  = 45106
$ws.Cells.Item(44, 5).Value  = 16
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100104
$ws.Cells.Item(44, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(44, 9).Value  = 100104003
$ws.Cells.Item(44, 10).Value = "Membrillo"
$ws.Cells.Item(44, 11).Value = "Champion"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 50
$ws.Cells.Item(44, 14).Value = 10000
$ws.Cells.Item(44, 15).Value = 10000
$ws.Cells.Item(44, 16).Value = 10000
$ws.Cells.Item(44, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(44, 18).Value = "Región del Maule"
$ws.Cells.Item(44, 19).Value = 556
$ws.Cells.Item(44, 20).Value = 18

# New row 45: "Especial" quality for the same week (serial 45104)
$ws.Cells.Item(45, 1).Value  = 7
$ws.Cells.Item(45, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(45, 3).Value  = "Ñuble"
$ws.Cells.Item(45, 4).Value  = 45104
$ws.Cells.Item(45, 5).Value  = 16
$ws.Cells.Item(45, 6).Value  = "Fruta"
$ws.Cells.Item(45, 7).Value  = 100104
$ws.Cells.Item(45, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(45, 9).Value  = 100104003
$ws.Cells.Item(45, 10).Value = "Membrillo"
$ws.Cells.Item(45, 11).Value = "Champion"
$ws.Cells.Item(45, 12).Value = "Especial"
$ws.Cells.Item(45, 13).Value = 80
$ws.Cells.Item(45, 14).Value = 12000
$ws.Cells.Item(45, 15).Value = 12000
$ws.Cells.Item(45, 16).Value = 12000
$ws.Cells.Item(45, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(45, 18).Value = "Región del Maule"
$ws.Cells.Item(45, 19).Value = 667
$ws.Cells.Item(45, 20).Value = 18

# Give the new date cells (column D) the same date formatting as the rest
# of the column.
$ws.Range("D44:D45").NumberFormat = $ws.Range("D46").NumberFormat

$ws.Range("A1").Select()
